$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.156.75'
$ws.Range('E2').Value = '  +2.17%  '

$ws.Range('D3').Value = '3.380.83'
$ws.Range('E3').Value = '  +1.63%  '

$ws.Range('E4').Value = '  +0.01%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '587.29'
$ws.Range('E5').Value = '  +1.01%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '180.39'
$ws.Range('E6').Value = '  +2.82%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.999'
$ws.Range('E7').Value = '  +0.05%  '

$ws.Range('E8').Value = '  +0.88%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.196'
$ws.Range('E9').Value = '  +8.38%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.588'
$ws.Range('E10').Value = '  +1.45%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '48.63'

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0000285'
$ws.Range('E12').Value = '  +5.04%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '685.83'
$ws.Range('E13').Value = '  -2.63%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '8.62'
$ws.Range('E14').Value = '  +2.18%  '

$ws.Range('D15').Value = '3.933.62'
$ws.Range('E15').Value = '  +1.54%  '

$ws.Range('D16').Value = '69.252.13'
$ws.Range('E16').Value = '  +2.26%  '

$ws.Range('B17').Value = 'WrappedEther'
$ws.Range('C17').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D17').Value = '3.398.54'
$ws.Range('E17').Value = '  +2.05%  '

$ws.Range('B18').Value = 'TRON'
$ws.Range('C18').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.120'
$ws.Range('E18').Value = '  +1.79%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '17.72'
$ws.Range('E19').Value = '  +2.04%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '11.36'
$ws.Range('E20').Value = '  +3.29%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.900'
$ws.Range('E21').Value = '  +0.83%  '

$ws.Range('E22').Value = '  +0.59%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '17.06'
$ws.Range('E23').Value = '  +0.68%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '104.62'
$ws.Range('E24').Value = '  +6.04%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '3.93'
$ws.Range('E25').Value = '  +1.34%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.72'
$ws.Range('E26').Value = '  +1.31%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.61'

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '34.31'
$ws.Range('E28').Value = '  +3.55%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '8.68'
$ws.Range('E29').Value = '  +1.80%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '6.96'
$ws.Range('E30').Value = '  -1.67%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '11.18'
$ws.Range('E31').Value = '  +1.88%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '556.77'
$ws.Range('E32').Value = '  -2.30%  '

$ws.Range('E33').Value = '  +9.74%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.106'
$ws.Range('E34').Value = '  +1.01%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '57.94'
$ws.Range('E35').Value = '  +0.87%  '

$ws.Range('E36').Value = '  +0.13%  '

$ws.Range('D37').Value = '3.701.00'
$ws.Range('E37').Value = '  -0.07%  '

$ws.Range('E38').Value = '  +7.07%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '34.83'
$ws.Range('E39').Value = '  +2.30%  '

$ws.Range('B40').Value = 'Stacks'
$ws.Range('C40').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.24'
$ws.Range('E40').Value = '  +1.40%  '

$ws.Range('B41').Value = 'PEPE'
$ws.Range('C41').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D41').Value = '0.0₃0704'
$ws.Range('E41').Value = '  +4.56%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.68'
$ws.Range('E42').Value = '  +1.84%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.339'
$ws.Range('E43').Value = '  +1.12%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0418'
$ws.Range('E44').Value = '  +3.06%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.27'
$ws.Range('E45').Value = '  -0.88%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.65'
$ws.Range('E46').Value = '  -1.12%  '

$ws.Range('E47').Value = '  +0.98%  '

$ws.Range('E48').Value = '  +5.18%  '

$ws.Range('E49').Value = '  -0.05%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '132.65'
$ws.Range('E50').Value = '  +2.70%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.57'
$ws.Range('E51').Value = '  -2.73%  '
